$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "0.9917") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.786.02"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.814.40"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "0.9917"
$ws.Range("E4").Value = "  -0.78%  "
$ws.Range("D5").Value = "241.44"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "0.6218"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "0.9928"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "0.07380"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "0.2902"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "22.75"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "0.07635"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "1.824.62"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "4.951"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "0.6601"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "82.33"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "0.000009482"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "5.980"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "28.802.17"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "12.45"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "222.32"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "0.9923"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "7.048"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").Value = "0.9935"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "158.22"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("D26").Value = "8.413"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "17.72"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "1.486"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "4.079"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "4.015"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "0.05401"
$ws.Range("E31").Value = "  +4.10%  "
$ws.Range("D32").Value = "1.185"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "1.832"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "0.7378"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "1.125"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "2.579"
$ws.Range("E36").Value = "  -4.39%  "
$ws.Range("D37").Value = "1.219.08"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").Value = "2.730"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "0.01767"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "6.608"
$ws.Range("E40").Value = "  +6.52%  "
$ws.Range("D41").Value = "0.8859"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "0.9924"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "100.81"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "64.65"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "0.5045"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "0.3998"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "8.868"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "0.07179"
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("D50").Value = "0.05769"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "1.641"
$ws.Range("E51").Value = "  +1.68%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so no residual number-format styling is left on the cells.
$ws.Range("D2:D51").Style = "Normal"
